# Add beam waist distance column
#
# Inserts a new column "r0" (beam waist distance) between the existing
# "w0" (column J) and "m" (old column K) columns, pushing the later
# parameter columns (m, n, D, T) one column to the right, and fills in
# the r0 values for the two rows (HEL1 = 1, HEL2 = 1000).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at K -- everything from K onward (m, n, D, T) shifts
# right to L, M, N, O. The new column inherits formatting from its left
# neighbour (J), which already carries the scientific-notation style used
# by the other numeric parameter columns.
$ws.Columns("K").EntireColumn.Insert()

# Header + data for the new "r0" column.
$ws.Range("K1").Value = "r0"
$ws.Range("K2").Value = 1
$ws.Range("K3").Value = 1000

# Size the new column to fit its contents, like the neighbouring columns.
$ws.Columns("K").AutoFit() | Out-Null

# Leave the selection where the new data was entered.
$ws.Range("K3").Select() | Out-Null
